$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53
$ws.Range("A$row").Value = "Squares of a Sorted Array"
$ws.Range("B$row").Value = "Array"
$ws.Range("C$row").Value = "No"
$ws.Range("D$row").Value = "Yes"
$ws.Range("E$row").Value = "Easy"
$ws.Range("F$row").Value = "Medium"
$ws.Range("G$row").Value = "977 - Squares of a Sorted Array"

$ws.Hyperlinks.Add($ws.Range("G$row"), "977 - Squares of a Sorted Array") | Out-Null
$ws.Range("G$row").Style = "Hyperlink"

$ws.Range("D2:G8,D9:F53").FormatConditions.Delete()

$cf = $ws.Range("D2:G8,D9:F53")

$ws.Range("E2:F53").Validation.Delete()
$ws.Range("E2:F53").Validation.Add(3, 1, 1, "Easy, Medium, Hard")

$ws.Range("C2:C53").Validation.Delete()
$ws.Range("C2:C53").Validation.Add(3, 1, 1, "Yes, No")

$ws.Range("D2:D53").Validation.Delete()
$ws.Range("D2:D53").Validation.Add(3, 1, 1, "Yes, No")

$ws.Range("B2:B53").Validation.Delete()
$ws.Range("B2:B53").Validation.Add(3, 1, 1, "Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design")
